$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- 1. Text changes (order matters: controls shared-string index order) ---
$ws.Range("B8").Value = "Mehr Matlab Funktionen?"
$ws.Range("D15").Value = "manche mehr"
$ws.Range("D18").Value = "Diskrete  Übtgsfkt. In Position Controller Guvp1"

$ws.Range("B3").ClearContents()
$ws.Range("B14").ClearContents()

# --- 2. Colour coding (order matters: controls fill/cellXf index order) ---
# Each group's formatting is created once on the first cell, then copy/pasted
# (format only) onto the remaining cells in the group so no transient,
# orphaned style records get minted along the way.

# 2a. Orange fill + orange font ("hidden text" marker cells): B5, B15
$ws.Range("B5").Interior.Color = 49407
$ws.Range("B5").Font.Color = 49407
$ws.Range("B5").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2b. Green fill: B3, B6, B7, B13, B14, B20
$ws.Range("B3").Interior.Color = 5287936
$ws.Range("B3").Copy()
foreach ($addr in @("B6","B7","B13","B14","B20")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# 2c. Red fill: B9, B10, B12, B17, B19, B21, B22
$ws.Range("B9").Interior.Color = 255
$ws.Range("B9").Copy()
foreach ($addr in @("B10","B12","B17","B19","B21","B22")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# 2d. Orange fill (plain, no special font): B4, B8, B11, B16, B18
$ws.Range("B4").Interior.Color = 49407
$ws.Range("B4").Copy()
foreach ($addr in @("B8","B11","B16","B18")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- 3. Column width ---
# Target stored width is 72.42578125 characters; the host quantizes
# ColumnWidth assignments to 1/6-character pixel steps, so 71.62 is the
# closest input that lands on the nearest achievable stored width (72.5).
$ws.Columns.Item(1).ColumnWidth = 71.62

# --- 4. Selection ---
$ws.Range("A5").Select()
